$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header "Label" in H1 (matches style of existing header cells)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New "Label" column values (0 = Control, 1 = MDD) for both iteration blocks
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1

# Refit: updated Prediction / Error / Cross Entropy Loss values (rows 2-11, block 1)
$ws.Range("D2").Value = 0.4499840218509325
$ws.Range("E2").Value = 0.4499840218509325

$ws.Range("D3").Value = 0.6538971239928671
$ws.Range("E3").Value = 0.6538971239928671

$ws.Range("D4").Value = 0.5406020981961104
$ws.Range("E4").Value = 0.5406020981961104

$ws.Range("D6").Value = 0.6285524556150874
$ws.Range("E6").Value = 0.6285524556150874

$ws.Range("D7").Value = 0.5165216104922983
$ws.Range("E7").Value = 0.4834783895077017

$ws.Range("D8").Value = 0.6056318138039664
$ws.Range("E8").Value = 0.3943681861960336

$ws.Range("D9").Value = 0.751006553934683
$ws.Range("E9").Value = 0.248993446065317

$ws.Range("D10").Value = 0.7097534151086355
$ws.Range("E10").Value = 0.2902465848913645

$ws.Range("F11").Value = 0.6888768672943115

# Refit: updated Cross Entropy Loss value (row 21, block 2)
$ws.Range("F21").Value = 0.6876498460769653
